# Update the raw experimental-data inputs on "alldata_1step" (columns E,
# rows 34-81, one row excluded at E61 which is unchanged in the source
# commit). Every other cell touched by this edit (F/J/K/L/N/O/P/Q columns,
# the scatter-chart caches, and the bar-chart caches) is formula-driven and
# recalculates automatically once the inputs change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("alldata_1step")

$ws.Range("E34").Value = 8.6601879999999998
$ws.Range("E35").Value = 9.0667150000000003
$ws.Range("E36").Value = 9.2164610000000007
$ws.Range("E37").Value = 9.1590299999999996

$ws.Range("E38").Value = 9.9775500000000008
$ws.Range("E39").Value = 9.4892839999999996
$ws.Range("E40").Value = 9.8151960000000003
$ws.Range("E41").Value = 9.8120130000000003

$ws.Range("E42").Value = 10.547658999999999
$ws.Range("E43").Value = 11.321642000000001
$ws.Range("E44").Value = 10.564242999999999
$ws.Range("E45").Value = 9.7783370000000005

$ws.Range("E46").Value = 11.832392
$ws.Range("E47").Value = 10.246069
$ws.Range("E48").Value = 9.4793640000000003
$ws.Range("E49").Value = 9.9483200000000007

$ws.Range("E50").Value = 8.7803909999999998
$ws.Range("E51").Value = 9.1737330000000004
$ws.Range("E52").Value = 9.2776730000000001
$ws.Range("E53").Value = 9.32578

$ws.Range("E54").Value = 10.041607000000001
$ws.Range("E55").Value = 9.8198120000000007
$ws.Range("E56").Value = 10.293556000000001
$ws.Range("E57").Value = 9.8548589999999994

$ws.Range("E58").Value = 10.363528000000001
$ws.Range("E59").Value = 9.9679210000000005
$ws.Range("E60").Value = 11.168996999999999
# E61 unchanged

$ws.Range("E62").Value = 11.357813
$ws.Range("E63").Value = 10.975279
$ws.Range("E64").Value = 11.072824000000001
$ws.Range("E65").Value = 11.075628999999999

$ws.Range("E66").Value = 11.109363
$ws.Range("E67").Value = 11.107533
$ws.Range("E68").Value = 11.141007
$ws.Range("E69").Value = 11.078530000000001

$ws.Range("E70").Value = 11.145773999999999
$ws.Range("E71").Value = 11.527557
$ws.Range("E72").Value = 11.008775999999999
$ws.Range("E73").Value = 12.356291000000001

$ws.Range("E74").Value = 11.589969
$ws.Range("E75").Value = 11.70552
$ws.Range("E76").Value = 11.39451
$ws.Range("E77").Value = 11.338858

$ws.Range("E78").Value = 12.432638000000001
$ws.Range("E79").Value = 12.816719000000001
$ws.Range("E80").Value = 12.255646
$ws.Range("E81").Value = 12.34074

# Reflect the author's final cursor position/selection on this sheet.
$ws.Activate()
$ws.Range("H101").Select()
